$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

$wsZh.Range("E3").Value = "2016-03-14 00:55:00"
$wsZh.Range("H3").Value = "2016-03-14 00:55:18"

$wsDe.Range("E3").Value = "2016-03-14 00:55:04"
$wsDe.Range("H3").Value = "2016-03-14 00:55:25"
